# Generate Report for Handback
# Refresh the handoff/handback timestamps for the "0b3b24c3-..." row
# (row 2) across the Overview / zh-cn / de-de sheets, as produced by a
# fresh handback-status report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" for the 0b3b24c3-... file
$wsOverview.Range("G2").Value = "2016-08-17 00:44:47"

# zh-cn: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the 0b3b24c3-... file
$wsZhCn.Range("H2").Value = "2016-08-17 00:44:42"
$wsZhCn.Range("K2").Value = "2016-08-17 00:44:59"

# de-de: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the 0b3b24c3-... file
$wsDeDe.Range("H2").Value = "2016-08-17 00:44:47"
$wsDeDe.Range("K2").Value = "2016-08-17 00:45:14"
